$d = $word.ActiveDocument

# The document's headers/footers each contain one inline picture (a logo).
# Commit swaps the cosmetic "name" the picture is known by:
#   - BTec logo  (in the headers) : image1.jpg -> image2.jpg
#   - Pearson logo (in the footers): image2.png -> image1.png
# (The embedded picture bytes / relationship targets are unchanged -
#  only the drawing's display name is being renamed.)

foreach ($sec in $d.Sections) {
    # Headers / Footers collections are fixed size 3:
    #   1 = wdHeaderFooterPrimary (first page) / wdHeaderFooterFirstPage depending on story,
    #   but regardless of which slot, iterate them all and only touch the
    #   ones that actually contain an inline picture.
    for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
        $hdr = $sec.Headers.Item($hi)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($si = 1; $si -le $shapes.Count; $si++) {
                $shapes.Item($si).Name = "image2.jpg"
            }
        }
    }

    for ($fi = 1; $fi -le $sec.Footers.Count; $fi++) {
        $ftr = $sec.Footers.Item($fi)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($si = 1; $si -le $shapes.Count; $si++) {
                $shapes.Item($si).Name = "image1.png"
            }
        }
    }
}
